$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 1234
$ws.Range("B5").Value = 11223344
$ws.Range("B6").Select()
